$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 268, shifting existing rows 268:350 down to 269:351
$ws.Rows.Item(268).Insert()

# Populate the newly inserted row 268 with the new price record
$ws.Range("A268").Value = 5
$ws.Range("B268").Value = "Macroferia Regional de Talca"
$ws.Range("C268").Value = "Maule"
$ws.Range("D268").Value = 44985
$ws.Range("E268").Value = 7
$ws.Range("F268").Value = "Fruta"
$ws.Range("G268").Value = 100108
$ws.Range("H268").Value = "Tropicales y subtropicales"
$ws.Range("I268").Value = 100108005
$ws.Range("J268").Value = "Piña"
$ws.Range("K268").Value = "Caramelo"
$ws.Range("L268").Value = "Segunda"
$ws.Range("M268").Value = 240
$ws.Range("N268").Value = 22000
$ws.Range("O268").Value = 22000
$ws.Range("P268").Value = 22000
$ws.Range("Q268").Value = "$/caja 14 unidades"
$ws.Range("R268").Value = "Ecuador"
$ws.Range("S268").Value = 1571
$ws.Range("T268").Value = 14
